# Swap the order of "System" and the email address in column G
# (Recorded By) wherever both are present as "System, <email>".
# After the edit these should read "<email>, System" instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.SpecialCells(11).Row   # xlCellTypeLastCell = 11
$col = 7                                    # column G

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2

    if ($val -ne $null -and $val -eq "System, dnasr281@gmail.com") {
        $cell.Value2 = "dnasr281@gmail.com, System"
    }
}
